$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# 1) Update the header row shading: 4EA72E -> 8DD873 (accent6, themeFillTint 99%)
$headerCell = $tbl.Cell(1, 1)
$headerCell.Shading.BackgroundPatternColor = 7592077   # 0x8DD873 -> RGB(141,216,115)

# 2) Remove the direct "not-bold" (b=0 / bCs=0) run/paragraph formatting on the
#    "PUESTO EN RANKING" number cells (rows 3-12, column 1). Re-insert a clean
#    paragraph with just the centered alignment and the plain text run.
$values = @("6", "10", "8", "4", "9", "1", "7", "5", "3", "2")
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 3 + $i
    $cell = $tbl.Cell($row, 1)
    $val = $values[$i]
    $xmlFrag = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>' + $val + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $cell.Range.InsertXML($xmlFrag)
}
